$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header E1
$ws.Range("E1").Value = "along"

# Row 2: All
$ws.Range("B2").Value = 0.0551760287523508
$ws.Range("C2").Value = 0.0387832854256928
$ws.Range("D2").Value = 0.0715687720790089
$ws.Range("E2").Value = "millionaire_tax_in_programTRUE"

# Row 3: Europe
$ws.Range("B3").Value = 0.0572694933776088
$ws.Range("C3").Value = 0.0338127770858666
$ws.Range("D3").Value = 0.0807262096693511
$ws.Range("E3").Value = "millionaire_tax_in_programTRUE"

# Row 4: France
$ws.Range("B4").Value = 0.00994252757416189
$ws.Range("C4").Value = -0.0443615343226938
$ws.Range("D4").Value = 0.0642465894710176
$ws.Range("E4").Value = "millionaire_tax_in_programTRUE"

# Row 5: Germany
$ws.Range("B5").Value = 0.0670965655915874
$ws.Range("C5").Value = 0.0172267698324775
$ws.Range("D5").Value = 0.116966361350697
$ws.Range("E5").Value = "millionaire_tax_in_programTRUE"

# Row 6: Italy
$ws.Range("B6").Value = 0.0817303877432258
$ws.Range("C6").Value = 0.0233570508655756
$ws.Range("D6").Value = 0.140103724620876
$ws.Range("E6").Value = "millionaire_tax_in_programTRUE"

# Row 7: Poland
$ws.Range("B7").Value = 0.0680533176453674
$ws.Range("C7").Value = -0.00248399960860055
$ws.Range("D7").Value = 0.138590634899335
$ws.Range("E7").Value = "millionaire_tax_in_programTRUE"

# Row 8: Spain
$ws.Range("B8").Value = 0.10158952222961
$ws.Range("C8").Value = 0.0337020950888191
$ws.Range("D8").Value = 0.169476949370402
$ws.Range("E8").Value = "millionaire_tax_in_programTRUE"

# Row 9: United Kingdom
$ws.Range("B9").Value = 0.0433125623099311
$ws.Range("C9").Value = -0.0150496713843895
$ws.Range("D9").Value = 0.101674796004252
$ws.Range("E9").Value = "millionaire_tax_in_programTRUE"

# Row 10: Switzerland
$ws.Range("B10").Value = -0.0145012916202922
$ws.Range("C10").Value = -0.0907107009404396
$ws.Range("D10").Value = 0.0617081176998553
$ws.Range("E10").Value = "millionaire_tax_in_programTRUE"

# Row 11: Japan
$ws.Range("B11").Value = 0.0342307930815836
$ws.Range("C11").Value = -0.000817353665763139
$ws.Range("D11").Value = 0.0692789398289303
$ws.Range("E11").Value = "millionaire_tax_in_programTRUE"

# Row 12: USA
$ws.Range("B12").Value = 0.062244897064752
$ws.Range("C12").Value = 0.033822594545063
$ws.Range("D12").Value = 0.090667199584441
$ws.Range("E12").Value = "millionaire_tax_in_programTRUE"
